$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New interview entries appended to the tracker (rows 23-24).
# Column B for row 24 is written first so that the shared-string table
# ends up with the same index ordering as the authored workbook.
$ws.Cells.Item(24, 2).Value = "Neosoft"
$ws.Cells.Item(23, 2).Value = "Nimbal, Andheri"
$ws.Cells.Item(23, 3).Value = "interface, sql queries"
$ws.Cells.Item(24, 3).Value = "optional.of vs optional.ofnullable, @springbootapplication, @transaction, singletone pattern, `ncode - [aaaa,bb,c] = a4b2c1 ."

$ws.Cells.Item(23, 1).Value = 45937
$ws.Cells.Item(24, 1).Value = 45942

# Row 24 wraps across two lines of text, matching the taller row used in the
# authored workbook.
$ws.Rows.Item(24).RowHeight = 43.2

# Reflect the saved selection/scroll position of the workbook after the new
# rows were entered.
$ws.Range("C25").Select() | Out-Null
